$wb = $excel.ActiveWorkbook

# --- "Set Values Here" sheet: update the government-revenue-accounting
# weights for the "carbon tax revenue" row (row 8). Deficit Spending (C)
# and Payroll Taxes (E) weights move from 0 to 5; the resulting row becomes
# Regular Spending=0, Deficit Spending=5, Household Taxes=5, Payroll Taxes=5,
# Corporate Taxes=0.
$wsSet = $wb.Worksheets.Item("Set Values Here")
[void]$wsSet.Activate()
$wsSet.Range("C8").Value = 5
$wsSet.Range("E8").Value = 5
[void]$wsSet.Range("C9").Select()

# --- "GRA-carbontax" sheet reads these weights via a TRANSPOSE array
# formula, so its cached values pick up the change automatically on
# recalculation. Just restore the on-sheet selection shown in the target.
$wsCarbon = $wb.Worksheets.Item("GRA-carbontax")
[void]$wsCarbon.Activate()
[void]$wsCarbon.Range("B5").Select()

# Restore the originally active sheet/tab.
[void]$wb.Worksheets.Item("About").Activate()

$wb.Save()
